$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 13's formatting (A13 carries the date number format) down to row 14
# so the new date cell reuses the existing style instead of creating a new one.
$ws.Range("A13").Copy($ws.Range("A14"))

# Populate the new row (row 14) with the data from the diff.
$ws.Range("A14").Value = 42619.890405092592
$ws.Range("B14").Value = 12
$ws.Range("C14").Value = 58
$ws.Range("D14").Value = 40
$ws.Range("E14").Value = 58
$ws.Range("F14").Value = 41
$ws.Range("G14").Value = 13296
$ws.Range("H14").Value = 28960
$ws.Range("I14").Value = 3273
$ws.Range("J14").Value = 487
$ws.Range("K14").Value = 339
$ws.Range("L14").Value = 54
$ws.Range("M14").Value = 39
$ws.Range("N14").Value = "Bag"
